# Refresh the cached "datetimeFigureOut" footer field (Insert > Header &
# Footer > Date and time) from 10/19/24 to 11/1/24 everywhere it is stored:
# the slide master and every slide layout each keep their own cached copy
# of the placeholder text.

$p = $ppt.ActivePresentation
$newDate = "11/1/24"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master.
$master = $p.Designs.Item(1).SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
